$wb = $excel.ActiveWorkbook

# "Croatia" is currently the last (active) sheet. Duplicate it to create
# the new "Greece" sheet, placing the copy right after "Croatia" - this
# carries over all formatting, column widths, merged cells and styles.
$croatia = $wb.Worksheets.Item("Croatia")
$croatia.Copy($null, $croatia)

# Right after Copy(), the newly created sheet is the active sheet.
$greece = $wb.ActiveSheet
$greece.Name = "Greece"

# Fill in the market-specific values for the new sheet.
$greece.Range("B2").Value = "Greece Market"
$greece.Range("B4").Value = "NGC-4119/T3167"

# "Croatia" is no longer the selected/active tab - select it and select
# all of its cells (matching the author deselecting cell B4 before
# saving), then re-activate "Greece" so it is the tab shown/selected
# when the workbook is opened.
$croatia.Select()
$croatia.Cells.Select()
$greece.Select()
